$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 365, shifting existing rows 365-479 down to 366-480.
$ws.Rows.Item(365).Insert()

# Populate the newly inserted row 365 with the new weekly price record.
$ws.Cells.Item(365, 1).Value = 7
$ws.Cells.Item(365, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(365, 3).Value = "Ñuble"
$ws.Cells.Item(365, 4).Value = 44841
$ws.Cells.Item(365, 5).Value = 16
$ws.Cells.Item(365, 6).Value = "Fruta"
$ws.Cells.Item(365, 7).Value = 100102
$ws.Cells.Item(365, 8).Value = "Cítricos"
$ws.Cells.Item(365, 9).Value = 100102005
$ws.Cells.Item(365, 10).Value = "Naranja"
$ws.Cells.Item(365, 11).Value = "Lane Late"
$ws.Cells.Item(365, 12).Value = "Primera"
$ws.Cells.Item(365, 13).Value = 160
$ws.Cells.Item(365, 14).Value = 7500
$ws.Cells.Item(365, 15).Value = 8000
$ws.Cells.Item(365, 16).Value = 7750
$ws.Cells.Item(365, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(365, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(365, 19).Value = 517
$ws.Cells.Item(365, 20).Value = 15
